$d = $word.ActiveDocument

# --- Step 1: 'Objetivos' paragraph (6) gets the old 'Programa resumido' text ---
$old1 = 'Orientar o/a aluno/aluna no início de sua trajetória universitária no curso de graduação em Engenharia na EEL/USP de modo que ele/ela seja capaz de a) identificar as oportunidades acadêmicas e as particularidades do seu curso; b) reconhecer, sob acompanhamento de um tutor, eventuais dificuldades ao longo do curso e compreender mecanismos para que estas sejam superadas, conduzindo o curso com o sucesso desejado; c) desenvolver habilidades técnicas e emocionais, ampliando as perspectivas de formação profissional por meio de atividades e encontros sistematizados.'
$new1 = 'Os cursos de engenharia, respectivos projetos pedagógicos e seus componentes curriculares, incluindo TCC e estágio obrigatório. Atividades extracurriculares. Identificação e aderência do/a aluno/aluna com o curso e com a profissão escolhida. O curso superior, a transição adolescente/jovem adulto e os desafios nos projetos de vida do/a aluno/aluna no início da graduação. Relação entre as disciplinas e o conhecimento a ser aplicado. Competências e habilidades desenvolvidas no seu curso de engenharia. Dimensões acadêmicas, socioculturais e científicas. Diversidade e inclusão. Organização dos estudos.'
$rng1 = $d.Paragraphs.Item(6).Range
$found1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 0, $false, $new1, 2)
Write-Output ("Step1 found=" + $found1)

# --- Step 2: insert the moved paragraphs' content before the docente-name list (paragraph 8) ---
$firstDocente = '144651 - Antonio Fernando Sartori'
$rng2 = $d.Paragraphs.Item(8).Range
$found2 = $rng2.Find.Execute($firstDocente, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
Write-Output ("Step2 found=" + $found2)
$insPoint = $d.Range($rng2.Start, $rng2.Start)
$insertBlock = ('Orientar o/a aluno/aluna no início de sua trajetória universitária no curso de graduação em Engenharia na EEL/USP de modo que ele/ela seja capaz de a) identificar as oportunidades acadêmicas e as particularidades do seu curso; b) reconhecer, sob acompanhamento de um tutor, eventuais dificuldades ao longo do curso e compreender mecanismos para que estas sejam superadas, conduzindo o curso com o sucesso desejado; c) desenvolver habilidades técnicas e emocionais, ampliando as perspectivas de formação profissional por meio de atividades e encontros sistematizados.' + [char]11 + 'Apresentação dos programas e serviços oferecidos pela USP voltados aos/às alunos/alunas e das oportunidades de realizar trabalhos extracurriculares. A dinâmica das aulas, ferramentas de interação. Desenvolvimento de atividades de grupo, com objetivo de desenvolver habilidades sócio-comportamentais através de colaboração em temas do curso relacionados à profissão escolhida. Áreas de atuação do curso de engenharia, competências e habilidades a serem desenvolvidas. Interdisciplinaridade e a relação entre as disciplinas e o conhecimento a ser aplicado. Planejamento de estudos. Formas de estudar e aprender.' + [char]11 + '(Descrever com clareza o processo de avaliação de aprendizagem para que o aluno compreenda como todos os elementos do plano de atividades se articulam e para que o professor possa realizar a gestão da aprendizagem na sua disciplina, com base em evidências do que o aluno aprendeu.) Método: Atividades realizadas na forma de dinâmicas de grupos, utilização de vídeos, textos, roda de discussão e/ou elaboração de painéis. Participação em encontros de orientação promovidos pelo Programa de Tutoria Acadêmica da EEL e a realização de atividades propostas pelo tutor/monitor/mentor incluindo trabalhos em equipe e estudos dirigidos.' + [char]11 + '(Descrever as experiências de aprendizagem ou como os conteúdos serão tratados do ponto de vista prático, a exemplo da realização de projetos, exercícios, seminários, visitas técnicas, simulações, atividade em laboratório, entre outras possibilidades) Participação ativa nos encontros, apresentação de estudos/pesquisa e de trabalhos realizados durante a disciplina, colaboração e engajamento nas atividades da disciplina. A nota final é dada pela média ponderada das notas obtidas nas diversas atividades propostas.' + [char]11 + 'Não se aplica' + [char]11 + 'A bibliografia será recomendada pelos docentes responsáveis e obtida na busca realizada pelos próprios alunos no início dos projetos. Seguem referências no tópico de mentoria.' + [char]11 + '' + [char]11 + '[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001.' + [char]11 + '[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27.' + [char]11 + '[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.' + [char]11 + '[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.' + [char]11 + '[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019.' + [char]11 + '')
$insPoint.InsertBefore($insertBlock)

# --- Step 3: remove the 6 docente names that moved elsewhere (Luiz..Sergio), from paragraph 8 ---
$hugo = '984972 - Hugo Ricardo Zschommler Sandim'
$rng3 = $d.Paragraphs.Item(8).Range
$found3 = $rng3.Find.Execute($hugo, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
Write-Output ("Step3 found=" + $found3)
$hugoEnd = $rng3.End
$docenteListPara = $d.Paragraphs.Item(8)
$paraEnd = $docenteListPara.Range.End
$delRange = $d.Range($hugoEnd, $paraEnd - 1)
Write-Output ("Step3 delRange=[" + $delRange.Text + "]")
$delRange.Text = ""

# --- Step 4: 'Programa resumido' paragraph (10) becomes Luiz's directory entry ---
$old4 = 'Os cursos de engenharia, respectivos projetos pedagógicos e seus componentes curriculares, incluindo TCC e estágio obrigatório. Atividades extracurriculares. Identificação e aderência do/a aluno/aluna com o curso e com a profissão escolhida. O curso superior, a transição adolescente/jovem adulto e os desafios nos projetos de vida do/a aluno/aluna no início da graduação. Relação entre as disciplinas e o conhecimento a ser aplicado. Competências e habilidades desenvolvidas no seu curso de engenharia. Dimensões acadêmicas, socioculturais e científicas. Diversidade e inclusão. Organização dos estudos.'
$new4 = '1176388 - Luiz Tadeu Fernandes Eleno'
$rng4 = $d.Paragraphs.Item(10).Range
$found4 = $rng4.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 0, $false, $new4, 2)
Write-Output ("Step4 found=" + $found4)

# --- Step 5: 'Programa' paragraph (12) becomes Maria's directory entry ---
$old5 = 'Apresentação dos programas e serviços oferecidos pela USP voltados aos/às alunos/alunas e das oportunidades de realizar trabalhos extracurriculares. A dinâmica das aulas, ferramentas de interação. Desenvolvimento de atividades de grupo, com objetivo de desenvolver habilidades sócio-comportamentais através de colaboração em temas do curso relacionados à profissão escolhida. Áreas de atuação do curso de engenharia, competências e habilidades a serem desenvolvidas. Interdisciplinaridade e a relação entre as disciplinas e o conhecimento a ser aplicado. Planejamento de estudos. Formas de estudar e aprender.'
$new5 = '7459752 - Maria Ismenia Sodero Toledo Faria'
$rng5 = $d.Paragraphs.Item(12).Range
$found5 = $rng5.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 0, $false, $new5, 2)
Write-Output ("Step5 found=" + $found5)

# --- Step 6: Avaliacao (14) - 'Metodo' answer becomes Miguel's directory entry ---
$old6 = '(Descrever com clareza o processo de avaliação de aprendizagem para que o aluno compreenda como todos os elementos do plano de atividades se articulam e para que o professor possa realizar a gestão da aprendizagem na sua disciplina, com base em evidências do que o aluno aprendeu.) Método: Atividades realizadas na forma de dinâmicas de grupos, utilização de vídeos, textos, roda de discussão e/ou elaboração de painéis. Participação em encontros de orientação promovidos pelo Programa de Tutoria Acadêmica da EEL e a realização de atividades propostas pelo tutor/monitor/mentor incluindo trabalhos em equipe e estudos dirigidos.'
$new6 = '5840622 - Miguel Justino Ribeiro Barboza'
$rng6 = $d.Paragraphs.Item(14).Range
$found6 = $rng6.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 0, $false, $new6, 2)
Write-Output ("Step6 found=" + $found6)

# --- Step 7: Avaliacao (14) - 'Criterio' answer becomes Sandra's directory entry ---
$old7 = '(Descrever as experiências de aprendizagem ou como os conteúdos serão tratados do ponto de vista prático, a exemplo da realização de projetos, exercícios, seminários, visitas técnicas, simulações, atividade em laboratório, entre outras possibilidades) Participação ativa nos encontros, apresentação de estudos/pesquisa e de trabalhos realizados durante a disciplina, colaboração e engajamento nas atividades da disciplina. A nota final é dada pela média ponderada das notas obtidas nas diversas atividades propostas.'
$new7 = '2166002 - Sandra Giacomin Schneider'
$rng7 = $d.Paragraphs.Item(14).Range
$found7 = $rng7.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 0, $false, $new7, 2)
Write-Output ("Step7 found=" + $found7)

# --- Step 8: Avaliacao (14) - 'Norma de recuperacao' answer becomes Sebastiao's directory entry ---
$old8 = 'Não se aplica'
$new8 = '1922320 - Sebastiao Ribeiro'
$rng8 = $d.Paragraphs.Item(14).Range
$found8 = $rng8.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 0, $false, $new8, 2)
Write-Output ("Step8 found=" + $found8)

# --- Step 9: 'Bibliografia' paragraph (16) becomes Sergio's directory entry ---
$old9 = ('A bibliografia será recomendada pelos docentes responsáveis e obtida na busca realizada pelos próprios alunos no início dos projetos. Seguem referências no tópico de mentoria.' + [char]11 + '' + [char]11 + '[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001.' + [char]11 + '[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27.' + [char]11 + '[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.' + [char]11 + '[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.' + [char]11 + '[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019.')
$new9 = '5840793 - Sérgio Schneider'
$rng9 = $d.Paragraphs.Item(16).Range
$found9 = $rng9.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 0, $false, $new9, 2)
Write-Output ("Step9 found=" + $found9)
